$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14: copy formatting from row 6 (same cell pattern: A,B,C,D,F,G,H) ---
$ws.Range("A6:D6").Copy()
$ws.Range("A14:D14").PasteSpecial(-4122)
$ws.Range("F6:H6").Copy()
$ws.Range("F14:H14").PasteSpecial(-4122)
$ws.Rows.Item(14).RowHeight = 90

# --- Row 15: copy formatting from row 13 (same cell pattern: A,B,C,D,F,G,H,I,J,K) ---
$ws.Range("A13:D13").Copy()
$ws.Range("A15:D15").PasteSpecial(-4122)
$ws.Range("F13:K13").Copy()
$ws.Range("F15:K15").PasteSpecial(-4122)
$ws.Rows.Item(15).RowHeight = 135

# --- Row 14 values (order matters for shared-string allocation) ---
$ws.Range("A14").Value = "Brought label"
$ws.Range("C14").Value = "1. Run the application`n2. Click Play`n3. Click on shop tab`n4. Click on one of the items in the shop`n5. Click the buy button"
$ws.Range("B14").Value = "Label pops up for 3 seconds then disapears to give feedback to the player"
$ws.Range("D14").Value = "The label pops up then disapears after 3 seconds"
$ws.Range("F14").Value = "Pass"
$ws.Range("G14").Value = "Caleb Gourley"
$ws.Range("H14").Value = 45435

# --- Row 15 values (order matters for shared-string allocation) ---
$ws.Range("B15").Value = "Updating money after tower is brought"
$ws.Range("A15").Value = "Updating money"
$ws.Range("C15").Value = "1. Run the application`n2. Click Play`n3. Click on shop tab`n4. Click on one of the items in the shop`n5. Click the buy button`n6. Click on game tab`n7. Check if money value has been updated"
$ws.Range("D15").Value = "The new money value is displayed when the item is brought"
$ws.Range("F15").Value = "Fail"
$ws.Range("G15").Value = "Caleb Gourley"
$ws.Range("H15").Value = 45435
$ws.Range("I15").Value = "Pass"
$ws.Range("J15").Value = "Caleb"
$ws.Range("K15").Value = 45435

# --- sheetView: selection / scroll position ---
$ws.Range("M19").Select()
$ws.Application.ActiveWindow.ScrollRow = 12

# --- workbook window size ---
$excel.Width = 28800
$excel.Height = 15555
$excel.Left = 4680
$excel.Top = 4680
